$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 80, shifting existing rows 80:141 down to 81:142.
$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80 with the new weekly record.
# Columns A, B, C, E-K are identical for every record in this sheet.
$ws.Cells.Item(80, 1).Value = 3
$ws.Cells.Item(80, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(80, 3).Value = "Coquimbo"
$ws.Cells.Item(80, 4).Value = 44818
$ws.Cells.Item(80, 5).Value = 5
$ws.Cells.Item(80, 6).Value = "Fruta"
$ws.Cells.Item(80, 7).Value = 100107
$ws.Cells.Item(80, 8).Value = "Otros"
$ws.Cells.Item(80, 9).Value = 100107011
$ws.Cells.Item(80, 10).Value = "Tuna"
$ws.Cells.Item(80, 11).Value = "Sin especificar"
$ws.Cells.Item(80, 12).Value = "Segunda"
$ws.Cells.Item(80, 13).Value = 54
$ws.Cells.Item(80, 14).Value = 25000
$ws.Cells.Item(80, 15).Value = 25000
$ws.Cells.Item(80, 16).Value = 25000
$ws.Cells.Item(80, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(80, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(80, 19).Value = 1562
$ws.Cells.Item(80, 20).Value = 16
